# Commit: swap the applied table style on the B1/B2 "types of financial
# documents" table (slide 5) from the deck's default custom table style
# to the built-in table style {07A92274-AE14-4BC2-8D18-0BAB65404BB2}.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$targetStyleId = "{07A92274-AE14-4BC2-8D18-0BAB65404BB2}"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle($targetStyleId)
    }
}
